# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 3 de Agosto de 2020 a las 18:52"

# Updated per-country statistics (row => column => new value)
$updates = @{
    4   = @{ B = 4827843; C = 14196; D = 2382805; E = 2286543; G = 130; H = 158495 }
    6   = @{ B = 1833922; C = 29220;                E = 595087;  G = 371; H = 38532 }
    11  = @{ B = 361493;  C = 1762;  D = 333976; E = 17810;  G = 99;  H = 9707 }
    12  = @{ B = 344134;  C = 3044;                                G = 9;   H = 28472 }
    20  = @{ B = 233851;  C = 995;   D = 217497; E = 10607;  G = 19;  H = 5747 }
    25  = @{ B = 117007;  C = 123;   D = 101578; E = 6482;   G = 2;   H = 8947 }
    46  = @{                         D = 47179;  E = 5845 }
    56  = @{ B = 36747;   C = 37;    D = 25669;  E = 9790;   G = 4;   H = 1288 }
    60  = @{ B = 31972;   C = 507;   D = 21901;  E = 8832;   G = 8;   H = 1239 }
    75  = @{ B = 16880;   C = 80;    D = 11706;  E = 4788;   G = 2;   H = 386 }
    95  = @{ B = 6864;    C = 9;                 E = 1554;   G = 1;   H = 118 }
    99  = @{ B = 5620;    C = 101;   D = 3031;   E = 2417;   G = 6;   H = 172 }
    102 = @{ B = 5062;    C = 177;               E = 3202;   G = 3;   H = 65 }
    104 = @{ B = 4737;    C = 75;                E = 3153;   G = 2;   H = 210 }
    141 = @{ B = 1218;    C = 5;     D = 1131;   E = 76 }
    146 = @{ B = 1155;    C = 5;                 E = 280 }
    169 = @{ B = 355;     C = 2;     D = 302;    E = 47 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
